$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=TO, F=SHORT, G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, L=DIFFERENCE (K and M stay blank)
$rows = @(
  @{ Row=360; A=359; B='Monday, Jan 16'; C='5:35 AM'; D='LO3828'; E='Warsaw'; F='(WAW)'; G='LOT '; H='E195'; I='(SP-LNK)'; J='5:47 AM'; L='0 hours, 12 minutes' },
  @{ Row=361; A=360; B='Monday, Jan 16'; C='6:00 AM'; D='KL1916'; E='Amsterdam'; F='(AMS)'; G='KLM '; H='E190'; I='(PH-EZR)'; J='6:05 AM'; L='0 hours, 5 minutes' },
  @{ Row=362; A=361; B='Monday, Jan 16'; C='6:00 AM'; D='SK1758'; E='Copenhagen'; F='(CPH)'; G='SAS '; H='A20N'; I='(EI-SIC)'; J='6:09 AM'; L='0 hours, 9 minutes' },
  @{ Row=363; A=362; B='Monday, Jan 16'; C='6:20 AM'; D='FR3277'; E='Oslo'; F='(TRF)'; G='Ryanair '; H='B738'; I='(SP-RKM)'; J='6:37 AM'; L='0 hours, 17 minutes' },
  @{ Row=364; A=363; B='Monday, Jan 16'; C='6:20 AM'; D='W61641'; E='Eindhoven'; F='(EIN)'; G='Wizz Air '; H='A320'; I='(HA-LYH)'; J='6:25 AM'; L='0 hours, 5 minutes' },
  @{ Row=365; A=364; B='Monday, Jan 16'; C='6:30 AM'; D='FR6095'; E='Gothenburg'; F='(GOT)'; G='Ryanair '; H='B738'; I='(SP-RSO)'; J='6:50 AM'; L='0 hours, 20 minutes' },
  @{ Row=366; A=365; B='Monday, Jan 16'; C='6:30 AM'; D='FR8782'; E='Cork'; F='(ORK)'; G='Ryanair '; H='B738'; I='(SP-RKQ)'; J='6:34 AM'; L='0 hours, 4 minutes' },
  @{ Row=367; A=366; B='Monday, Jan 16'; C='6:30 AM'; D='W61731'; E='Stockholm'; F='(NYO)'; G='Wizz Air '; H='A321'; I='(HA-LXL)'; J='6:46 AM'; L='0 hours, 16 minutes' },
  @{ Row=368; A=367; B='Monday, Jan 16'; C='6:40 AM'; D='FR6844'; E='Copenhagen'; F='(CPH)'; G='Ryanair '; H='B738'; I='(SP-RSL)'; J='6:56 AM'; L='0 hours, 16 minutes' },
  @{ Row=369; A=368; B='Monday, Jan 16'; C='6:55 AM'; D='FR6111'; E='Lublin'; F='(LUZ)'; G='Ryanair '; H='B738'; I='(SP-RSW)'; J='6:53 AM'; L='0 hours, -2 minutes' },
  @{ Row=370; A=369; B='Monday, Jan 16'; C='6:55 AM'; D='LH1379'; E='Frankfurt'; F='(FRA)'; G='Lufthansa '; H='CRJ9'; I='(D-ACNJ)'; J='7:01 AM'; L='0 hours, 6 minutes' },
  @{ Row=371; A=370; B='Monday, Jan 16'; C='7:15 AM'; D='E45647'; E='Hurghada'; F='(HRG)'; G='Enter Air '; H='B738'; I='(SP-ESD)'; J='7:44 AM'; L='0 hours, 29 minutes' },
  @{ Row=372; A=371; B='Monday, Jan 16'; C='7:50 AM'; D='LH1647'; E='Munich'; F='(MUC)'; G='Lufthansa '; H='CRJ9'; I='(D-ACKF)'; J='7:51 AM'; L='0 hours, 1 minutes' },
  @{ Row=373; A=372; B='Monday, Jan 16'; C='7:55 AM'; D='FR4106'; E='Wroclaw'; F='(WRO)'; G='Ryanair '; H='B738'; I='(SP-RKI)'; J='8:11 AM'; L='0 hours, 16 minutes' },
  @{ Row=374; A=373; B='Monday, Jan 16'; C='7:55 AM'; D='W61745'; E='Bergen'; F='(BGO)'; G='Wizz Air '; H='A21N'; I='(9H-WBU)'; J='7:59 AM'; L='0 hours, 4 minutes' },
  @{ Row=375; A=374; B='Monday, Jan 16'; C='8:50 AM'; D='LO3832'; E='Warsaw'; F='(WAW)'; G='LOT '; H='E190'; I='(SP-LME)'; J='9:03 AM'; L='0 hours, 13 minutes' },
  @{ Row=376; A=375; B='Monday, Jan 16'; C='9:55 AM'; D='FR6113'; E='London'; F='(STN)'; G='Ryanair '; H='B738'; I='(SP-RSW)'; J='9:58 AM'; L='0 hours, 3 minutes' },
  @{ Row=377; A=376; B='Monday, Jan 16'; C='10:10 AM'; D='FR6123'; E='Edinburgh'; F='(EDI)'; G='Ryanair '; H='B738'; I='(SP-RSL)'; J='10:18 AM'; L='0 hours, 8 minutes' }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 1).Value = $r.A
  $ws.Cells.Item($r.Row, 2).Value = $r.B
  $ws.Cells.Item($r.Row, 3).Value = $r.C
  $ws.Cells.Item($r.Row, 4).Value = $r.D
  $ws.Cells.Item($r.Row, 5).Value = $r.E
  $ws.Cells.Item($r.Row, 6).Value = $r.F
  $ws.Cells.Item($r.Row, 7).Value = $r.G
  $ws.Cells.Item($r.Row, 8).Value = $r.H
  $ws.Cells.Item($r.Row, 9).Value = $r.I
  $ws.Cells.Item($r.Row, 10).Value = $r.J
  $ws.Cells.Item($r.Row, 11).Interior.Pattern = -4142
  $ws.Cells.Item($r.Row, 12).Value = $r.L
  $ws.Cells.Item($r.Row, 13).Interior.Pattern = -4142
}
